# Updated RAD test data
# Re-applies the FEINandFEINSSNlessThan9Error RAD test-data refresh:
#   - the Date column (B) for the existing 47 result rows (rows 2-48) is
#     rewritten with the timestamps from the re-run executed on
#     2025-02-11 19:48-19:57 EST
#   - five new result rows (49-53) are appended for the new
#     "Digital Advertising Gross Revenues" / "PTE Composite" tax types

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New timestamps for the Date column of the existing rows (2..48), in
# the same order the original log rows were produced.
# ---------------------------------------------------------------------
$newDates = @(
    "Tue Feb 11 19:48:59 EST 2025",
    "Tue Feb 11 19:49:09 EST 2025",
    "Tue Feb 11 19:49:18 EST 2025",
    "Tue Feb 11 19:49:28 EST 2025",
    "Tue Feb 11 19:49:38 EST 2025",
    "Tue Feb 11 19:49:48 EST 2025",
    "Tue Feb 11 19:49:57 EST 2025",
    "Tue Feb 11 19:50:07 EST 2025",
    "Tue Feb 11 19:50:17 EST 2025",
    "Tue Feb 11 19:50:26 EST 2025",
    "Tue Feb 11 19:50:36 EST 2025",
    "Tue Feb 11 19:50:45 EST 2025",
    "Tue Feb 11 19:50:55 EST 2025",
    "Tue Feb 11 19:51:05 EST 2025",
    "Tue Feb 11 19:51:14 EST 2025",
    "Tue Feb 11 19:51:24 EST 2025",
    "Tue Feb 11 19:51:34 EST 2025",
    "Tue Feb 11 19:51:44 EST 2025",
    "Tue Feb 11 19:51:54 EST 2025",
    "Tue Feb 11 19:52:05 EST 2025",
    "Tue Feb 11 19:52:14 EST 2025",
    "Tue Feb 11 19:52:24 EST 2025",
    "Tue Feb 11 19:52:34 EST 2025",
    "Tue Feb 11 19:52:44 EST 2025",
    "Tue Feb 11 19:52:54 EST 2025",
    "Tue Feb 11 19:53:04 EST 2025",
    "Tue Feb 11 19:53:14 EST 2025",
    "Tue Feb 11 19:53:23 EST 2025",
    "Tue Feb 11 19:53:33 EST 2025",
    "Tue Feb 11 19:53:43 EST 2025",
    "Tue Feb 11 19:53:52 EST 2025",
    "Tue Feb 11 19:54:02 EST 2025",
    "Tue Feb 11 19:54:12 EST 2025",
    "Tue Feb 11 19:54:21 EST 2025",
    "Tue Feb 11 19:54:31 EST 2025",
    "Tue Feb 11 19:54:41 EST 2025",
    "Tue Feb 11 19:54:50 EST 2025",
    "Tue Feb 11 19:55:00 EST 2025",
    "Tue Feb 11 19:55:10 EST 2025",
    "Tue Feb 11 19:55:19 EST 2025",
    "Tue Feb 11 19:55:29 EST 2025",
    "Tue Feb 11 19:55:39 EST 2025",
    "Tue Feb 11 19:55:48 EST 2025",
    "Tue Feb 11 19:55:58 EST 2025",
    "Tue Feb 11 19:56:08 EST 2025",
    "Tue Feb 11 19:56:18 EST 2025",
    "Tue Feb 11 19:56:27 EST 2025",
    "Tue Feb 11 19:56:37 EST 2025",
    "Tue Feb 11 19:56:47 EST 2025",
    "Tue Feb 11 19:56:57 EST 2025",
    "Tue Feb 11 19:57:06 EST 2025",
    "Tue Feb 11 19:57:16 EST 2025"
)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newDates[$i]
}

# ---------------------------------------------------------------------
# Append the five new rows (49-53) produced by the new test steps.
# Columns A/B keep the workbook's default (unbordered) style, while
# C:F reuse the bordered/wrap-text style ("style 1") used by every
# other data row, matching the existing sheet layout.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row=49; A="Pass"; B="Tue Feb 11 19:56:37 EST 2025"; C="Y"; D="Existing Liability with Notice/Invoice Number"; E="Digital Advertising Gross Revenues"; F=$null },
    @{ Row=50; A="Pass"; B="Tue Feb 11 19:56:47 EST 2025"; C="Y"; D="New Tax Return Amount Due";                     E="Digital Advertising Gross Revenues"; F=$null },
    @{ Row=51; A="Pass"; B="Tue Feb 11 19:56:57 EST 2025"; C="Y"; D="Existing Liability with Notice/Invoice Number"; E="PTE Composite";                     F="Y" },
    @{ Row=52; A="Pass"; B="Tue Feb 11 19:57:06 EST 2025"; C="Y"; D="New Tax Return Amount Due";                     E="IFTA Tax";                         F="Y" },
    @{ Row=53; A="Pass"; B="Tue Feb 11 19:57:16 EST 2025"; C="Y"; D="New Tax Return Amount Due";                     E="PTE Composite";                     F="Y" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Reset A:B to the default (no border) style before writing, since a
    # brand new row would otherwise inherit the bordered column style.
    $ws.Range("A" + $row + ":B" + $row).Style = "Normal"

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B

    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    if ($r.F -ne $null) {
        $ws.Cells.Item($row, 6).Value = $r.F
    }
}

# Move the selection down to the newly-appended rows, matching the
# on-disk sheetView state captured after the refreshed RAD run.
$ws.Range("C43:C53").Select() | Out-Null
